# "update load file excel"
#
# - Sheet2 rows 11..18 (A:F) get filled in with a copy of Sheet1!A2:F9
#   (the "Ngày/Nhật" kanji-study block that was already present on Sheet1).
# - Sheet2!C10:E10 get three brand-new strings (とおく / 遠く / " xa (adv)").
# - The active sheet/tab flips from Sheet1 to Sheet2, with new selections
#   and scroll positions on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Fill in Sheet2 rows 11-18 with a copy of Sheet1's A2:F9 block ---
$ws1.Range("A2:F9").Copy()
$ws2.Range("A11").PasteSpecial()

# --- New vocabulary row inserted at Sheet2!C10:E10 ---
$ws2.Range("C10").Value = "とおく"
$ws2.Range("D10").Value = "遠く"
$ws2.Range("E10").Value = " xa (adv)"

# --- View/selection state ---
# Sheet1 loses the tab-selected flag, scrolls up, and selects A7:F9.
$ws1.Activate() | Out-Null
$ws1.Range("A7:F9").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 511
$excel.ActiveWindow.ScrollColumn = 1

# Sheet2 becomes the active/selected tab, with A16 selected.
$ws2.Activate() | Out-Null
$ws2.Range("A16").Select() | Out-Null
